# Add two new client rows (21 and 22) to the "Clientes" sheet, matching
# the pattern already used by the existing rows (columns A,B,C,E,F,G,H are
# populated; column D "Endereço" is left empty for these records).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data to append: row number -> ordered values for columns A,B,C,E,F,G,H
$newRows = @(
    @{ Row = 21; A = "BRUNO DE FRAGA"; B = "123123";   C = "123123";   E = "92320-195"; F = "12312312@123123"; G = "123123";   H = "Rua 3 Pinheiros I, 27" },
    @{ Row = 22; A = "BRUNO DE FRAGA"; B = "12312312"; C = "12312312"; E = "92320-195"; F = "joanues@gmail.com"; G = "12312312"; H = "Rua 3 Pinheiros I, 27" }
)

foreach ($entry in $newRows) {
    $r = $entry.Row

    # Force text storage (matches the existing "number stored as text" data
    # already present throughout the sheet) before writing the values, so
    # purely numeric-looking strings like "123123" are kept as text and not
    # coerced into numbers.
    foreach ($col in 1,2,3,5,6,7,8) {
        $ws.Cells.Item($r, $col).NumberFormat = "@"
    }

    $ws.Cells.Item($r, 1).Value = $entry.A
    $ws.Cells.Item($r, 2).Value = $entry.B
    $ws.Cells.Item($r, 3).Value = $entry.C
    $ws.Cells.Item($r, 5).Value = $entry.E
    $ws.Cells.Item($r, 6).Value = $entry.F
    $ws.Cells.Item($r, 7).Value = $entry.G
    $ws.Cells.Item($r, 8).Value = $entry.H
}

Write-Host "New used range: $($ws.UsedRange.Address())"
